$wb = $excel.ActiveWorkbook

# --- OFF sheet (row 2: "H") ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 195
$wsOff.Range("C2").Value = 123
$wsOff.Range("D2").Value = 48
$wsOff.Range("E2").Value = 24
$wsOff.Range("G2").Value = 3

# --- DEF sheet (row 2: "H") ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 165
$wsDef.Range("C2").Value = 114
$wsDef.Range("F2").Value = 3
